$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new value would otherwise be
# auto-detected as a number by Excel, so they stay text like the rest of column D.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.365.02"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "3.091.42"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "386.57"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "103.47"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("D10").Value = "36.90"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "0.0856"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "3.579.07"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").Value = "18.50"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").Value = "7.77"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "3.096.13"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("D17").Value = "0.990"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "10.64"
$ws.Range("D19").Value = "51.434.06"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "3.22"
$ws.Range("E20").Value = "  +4.60%  "
$ws.Range("D21").Value = "12.48"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "70.11"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "265.56"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").Value = "8.00"
$ws.Range("E26").Value = "  -3.09%  "
$ws.Range("D27").Value = "27.33"
$ws.Range("E27").Value = "  +4.00%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "7.19"
$ws.Range("E28").Value = "  -5.98%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -5.09%  "
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").Value = "10.40"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").Value = "35.60"
$ws.Range("E33").Value = "  +4.32%  "
$ws.Range("D34").Value = "0.0474"
$ws.Range("E34").Value = "  +5.59%  "
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").Value = "50.00"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("D39").Value = "0.291"
$ws.Range("E39").Value = "  -2.26%  "
$ws.Range("D40").Value = "129.50"
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "16.51"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").Value = "3.82"
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("D46").Value = "21.98"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "2.51"
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("D48").Value = "2.09"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "2.072.92"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.928"
$ws.Range("E50").Value = "  +18.06%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "0.0330"
$ws.Range("E51").Value = "  +2.80%  "
